# Add a new SKU value (10126309) as the next row in the vp_sku_list sheet,
# matching the formatting (Roboto, 11pt, color #202124) used for this entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A86").Value = 10126309
$ws.Range("A86").Font.Name = "Roboto"
$ws.Range("A86").Font.Color = 0x20 + (0x21 * 256) + (0x24 * 65536)  # RGB(0x20,0x21,0x24) -> FF202124
[void]$ws.Range("A86").Select()
